$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B18").Value = 14
$ws.Range("C18").Value = 1515
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = "lower long"
$ws.Range("F18").Value = "Extrusion"
$ws.Range("G18").Value = 48.19
$ws.Range("H18").Value = "7010 cut to length. 7040 sides A and C, both ends."

$ws.Range("H18").Select()
